# Round 2 data updates: update match dates in column A
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 45738
    3  = 45739
    4  = 45736
    5  = 45737
    6  = 45738
    7  = 45739
    9  = 45738
    11 = 45736
    12 = 45739
    13 = 45739
    14 = 45738
    15 = 45738
    16 = 45738
    17 = 45739
    18 = 45739
    19 = 45737
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
